$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-30 Friday" "2026-01-31 Saturday"

Replace-Text "360÷5=72, 0" "589÷2=294, 1"
Replace-Text "995÷8=124, 3" "568÷8=71, 0"
Replace-Text "936÷6=156, 0" "412÷7=58, 6"
Replace-Text "807÷9=89, 6" "792÷4=198, 0"
Replace-Text "864÷6=144, 0" "711÷2=355, 1"

Replace-Text "352÷6=58, 4" "281÷3=93, 2"
Replace-Text "949÷4=237, 1" "736÷9=81, 7"
Replace-Text "643÷6=107, 1" "105÷4=26, 1"
Replace-Text "158÷2=79, 0" "118÷3=39, 1"
Replace-Text "265÷4=66, 1" "507÷9=56, 3"

Replace-Text "511÷6=85, 1" "698÷8=87, 2"
Replace-Text "424÷6=70, 4" "827÷6=137, 5"
Replace-Text "538÷8=67, 2" "444÷4=111, 0"
Replace-Text "523÷5=104, 3" "238÷7=34, 0"
Replace-Text "791÷2=395, 1" "221÷6=36, 5"

Replace-Text "659÷6=109, 5" "904÷6=150, 4"
Replace-Text "366÷7=52, 2" "908÷5=181, 3"
Replace-Text "310÷4=77, 2" "467÷5=93, 2"
Replace-Text "840÷7=120, 0" "277÷9=30, 7"
Replace-Text "636÷6=106, 0" "774÷5=154, 4"

Replace-Text "726÷2=363, 0" "645÷7=92, 1"
Replace-Text "121÷3=40, 1" "842÷8=105, 2"
Replace-Text "113÷3=37, 2" "526÷8=65, 6"
Replace-Text "163÷4=40, 3" "580÷4=145, 0"
Replace-Text "449÷3=149, 2" "440÷4=110, 0"

Write-Output "Done replacing all text."
